# Update search in DG
# - filterCurrentList(listType, predicate) -> filterCurrentList(listType, p)
# - setPredicate(predicate)                -> setPredicate(p)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape "TextBox 67" (id 123): filterCurrentList(listType, predicate)
$shape1 = $s.Shapes.Item(10)
$tr1 = $shape1.TextFrame.TextRange
$old1 = $tr1.Text
$needle1 = ", predicate)"
$start1 = $old1.IndexOf($needle1) + 1
$sub1 = $tr1.Characters($start1, $needle1.Length)
$sub1.Text = ", p)"

# Shape "TextBox 67" (id 146): setPredicate(predicate)
$shape2 = $s.Shapes.Item(15)
$tr2 = $shape2.TextFrame.TextRange
$old2 = $tr2.Text
$needle2 = "(predicate)"
$start2 = $old2.IndexOf($needle2) + 1
$sub2 = $tr2.Characters($start2, $needle2.Length)
$sub2.Text = "(p)"
